$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "themes" sheet: the t1/level-1 row's minseconds value (column I)
#    changes from 18 to 8 -- part of the new minbeats/minseconds support.
# ---------------------------------------------------------------------
$themes = $wb.Worksheets.Item("themes")
$themes.Cells.Item(2, 9).Value = 8

# ---------------------------------------------------------------------
# 2) "regions" sheet: a new region row "test-t1-1" (group "test", theme
#    t1, level 1) is inserted above the existing "test-t1-2" row,
#    pushing the rest of that block (test-t1-2/-3, nothing, test-t4-1/-2)
#    down by one row.
# ---------------------------------------------------------------------
$regions = $wb.Worksheets.Item("regions")

$regions.Rows.Item(12).Insert()
$regions.Rows.Item(12).Clear()

# Copy the formatting of the row directly below (the shifted former row
# 12) onto the matching cells of the new row, so styles line up exactly
# with the rest of this block.
$regions.Cells.Item(13, 1).Copy()
$regions.Cells.Item(12, 1).PasteSpecial(-4122)   # xlPasteFormats
$regions.Cells.Item(13, 3).Copy()
$regions.Cells.Item(12, 3).PasteSpecial(-4122)
$regions.Cells.Item(13, 4).Copy()
$regions.Cells.Item(12, 4).PasteSpecial(-4122)
$regions.Cells.Item(13, 9).Copy()
$regions.Cells.Item(12, 9).PasteSpecial(-4122)
$regions.Cells.Item(13, 13).Copy()
$regions.Cells.Item(12, 13).PasteSpecial(-4122)
$regions.Cells.Item(13, 14).Copy()
$regions.Cells.Item(12, 14).PasteSpecial(-4122)
$regions.Cells.Item(13, 15).Copy()
$regions.Cells.Item(12, 15).PasteSpecial(-4122)

$regions.Cells.Item(12, 1).Value = "test"         # A12 group
$regions.Cells.Item(12, 3).Value = "test-t1-1"    # C12 region
$regions.Cells.Item(12, 4).Value = "n"             # D12 gps
$regions.Cells.Item(12, 9).Value = 1               # I12 priority
$regions.Cells.Item(12, 13).Value = "nothing"      # M12 disable
$regions.Cells.Item(12, 14).Value = "t1"           # N12 theme
$regions.Cells.Item(12, 15).Value = 1              # O12 level

# Restore the selection to where the author last left it.
$regions.Range("P12").Select() | Out-Null
